# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.197.65"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "1.911.03"
$ws.Range("E3").Value = "  +2.03%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'314.16"
$ws.Range("E5").Value = "  +0.88%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").Value = "'0.5063"
$ws.Range("E7").Value = "  +0.45%  "

$ws.Range("D8").Value = "'0.3923"
$ws.Range("E8").Value = "  +0.53%  "

$ws.Range("D9").Value = "'0.09280"
$ws.Range("E9").Value = "  -2.44%  "

$ws.Range("D10").Value = "'1.141"
$ws.Range("E10").Value = "  +0.04%  "

$ws.Range("E11").Value = "  +2.92%  "

$ws.Range("D12").Value = "'6.408"
$ws.Range("E12").Value = "  -0.57%  "

$ws.Range("D13").Value = "'20.98"
$ws.Range("E13").Value = "  +0.16%  "

$ws.Range("D14").Value = "1.905.28"
$ws.Range("E14").Value = "  +1.86%  "

$ws.Range("D15").Value = "'7.329"
$ws.Range("E15").Value = "  -0.75%  "

$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  -0.17%  "

$ws.Range("D17").Value = "'0.00001123"
$ws.Range("E17").Value = "  -0.07%  "

$ws.Range("D18").Value = "'92.65"
$ws.Range("E18").Value = "  +0.05%  "

$ws.Range("E19").Value = "  +0.48%  "

$ws.Range("D20").Value = "'18.01"
$ws.Range("E20").Value = "  +1.92%  "

$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("D22").Value = "'6.237"
$ws.Range("E22").Value = "  +1.20%  "

$ws.Range("D23").Value = "28.270.19"
$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").Value = "'11.54"
$ws.Range("E24").Value = "  +2.21%  "

$ws.Range("D25").Value = "'2.328"
$ws.Range("E25").Value = "  +2.50%  "

$ws.Range("D26").Value = "'2.591"
$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("D27").Value = "2.141.82"
$ws.Range("E27").Value = "  +2.66%  "

$ws.Range("D28").Value = "'21.14"
$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").Value = "'158.21"
$ws.Range("E29").Value = "  -0.48%  "

$ws.Range("D30").Value = "'127.35"
$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("D31").Value = "'1.106"
$ws.Range("E31").Value = "  +3.66%  "

$ws.Range("D32").Value = "'0.1074"
$ws.Range("E32").Value = "  +1.01%  "

$ws.Range("D33").Value = "'5.675"
$ws.Range("E33").Value = "  +0.93%  "

$ws.Range("D34").Value = "'3.611"
$ws.Range("E34").Value = "  -0.32%  "

$ws.Range("D35").Value = "'9.699"
$ws.Range("E35").Value = "  +1.95%  "

$ws.Range("D36").Value = "'0.06702"
$ws.Range("E36").Value = "  -0.65%  "

$ws.Range("D37").Value = "'0.02440"
$ws.Range("E37").Value = "  +0.92%  "

$ws.Range("D38").Value = "'0.2230"
$ws.Range("E38").Value = "  +2.20%  "

$ws.Range("D39").Value = "'1.245"
$ws.Range("E39").Value = "  +0.85%  "

$ws.Range("D40").Value = "'1.281"
$ws.Range("E40").Value = "  +8.12%  "

$ws.Range("D41").Value = "'0.6576"
$ws.Range("E41").Value = "  +3.49%  "

$ws.Range("D42").Value = "'11.56"

$ws.Range("D43").Value = "'5.017"
$ws.Range("E43").Value = "  +0.41%  "

$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("D45").Value = "'0.6157"
$ws.Range("E45").Value = "  +2.91%  "

$ws.Range("D46").Value = "'13.37"
$ws.Range("E46").Value = "  -1.21%  "

$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.725"
$ws.Range("E47").Value = "  +1.90%  "

$ws.Range("B48").Value = "WEMIXTOKEN"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "'1.296"
$ws.Range("E48").Value = "  +1.50%  "

$ws.Range("D49").Value = "'2.030"
$ws.Range("E49").Value = "  +1.73%  "

$ws.Range("D50").Value = "'122.20"
$ws.Range("E50").Value = "  -0.93%  "

$ws.Range("D51").Value = "'1.190"
$ws.Range("E51").Value = "  -0.44%  "

